$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Update shared-string text values
$ws.Range("B3").Value = "Отбор пациентов на ВМП"
$ws.Range("B4").Value = "Консилиум"

# Update date/numeric values in column D (rows 2-4)
$ws.Range("D2").Value = 2958465
$ws.Range("D3").Value = 2958465
$ws.Range("D4").Value = 2958465
